$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (L1) - update raw inputs; dependent formulas (H4, K4) recalc automatically
$ws.Range("D4").Value = 1580
$ws.Range("G4").Value = 2050

# Row 5 (L2) - update raw inputs; dependent formulas (H5, I5, K5, L5) recalc automatically
$ws.Range("C5").Value = 1570
$ws.Range("D5").Value = 1550

# Row 6 (L3) - update raw inputs; dependent formulas (H6, I6, K6, L6) recalc automatically
$ws.Range("C6").Value = 1350
$ws.Range("G6").Value = 2150

# Row 7 (R1) - update raw inputs; dependent formulas (H7, I7, K7, L7) recalc automatically
$ws.Range("C7").Value = 1410
$ws.Range("D7").Value = 1380

# Row 8 (R2) - update raw inputs; dependent formulas (H8, I8, K8, L8) recalc automatically
$ws.Range("C8").Value = 1600
$ws.Range("D8").Value = 1250

# Update the saved selection/active cell to H12 (sheet view no longer scrolled to E1)
$ws.Range("H12").Select()
